$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("19Tto20TMap")

$ws.Range("A2").Value = 'GEND_GBV.PE.T'
$ws.Range("A3").Value = 'GEND_GBV.S.T'
$ws.Range("A4").Value = 'HTS_INDEX_COM.New.Neg.T'
$ws.Range("A5").Value = 'HTS_INDEX_COM.New.Pos.T'
$ws.Range("A6").Value = 'HTS_RECENT.IndexCom.T'
$ws.Range("A7").Value = 'HTS_INDEX_FAC.New.Neg.T'
$ws.Range("A8").Value = 'HTS_INDEX_FAC.New.Pos.T'
$ws.Range("A9").Value = 'HTS_RECENT.IndexFac.T'
$ws.Range("A10").Value = 'HTS_SELF.T'
$ws.Range("A11").Value = 'HTS_TST.EW.Neg.T'
$ws.Range("A12").Value = 'HTS_TST.EW.Pos.T'
$ws.Range("A13").Value = 'HTS_RECENT.EW.T'
$ws.Range("A14").Value = 'HTS_TST.Inpat.Neg.T'
$ws.Range("A15").Value = 'HTS_TST.Inpat.Pos.T'
$ws.Range("A16").Value = 'HTS_RECENT.Inpat.T'
$ws.Range("A17").Value = 'HTS_TST.Maln.Neg.T'
$ws.Range("A18").Value = 'HTS_TST.Maln.Pos.T'
$ws.Range("A19").Value = 'HTS_TST.MobileCom.Neg.T'
$ws.Range("A20").Value = 'HTS_TST.MobileCom.Pos.T'
$ws.Range("A21").Value = 'HTS_RECENT.MobileCom.T'
$ws.Range("A22").Value = 'HTS_TST.OtherCom.Neg.T'
$ws.Range("A23").Value = 'HTS_TST.OtherCom.Pos.T'
$ws.Range("A24").Value = 'HTS_RECENT.OtherCom.T'
$ws.Range("A25").Value = 'HTS_TST.Other.Neg.T'
$ws.Range("A26").Value = 'HTS_TST.Other.Pos.T'
$ws.Range("A27").Value = 'HTS_RECENT.Other.T'
$ws.Range("A28").Value = 'HTS_TST.Peds.Neg.T'
$ws.Range("A29").Value = 'HTS_TST.Peds.Pos.T'
$ws.Range("A30").Value = 'HTS_TST.PostANC1.Neg.T'
$ws.Range("A31").Value = 'HTS_TST.PostANC1.Pos.T'
$ws.Range("A32").Value = 'HTS_RECENT.PostANC1.T'
$ws.Range("A33").Value = 'HTS_TST.STI.Neg.T'
$ws.Range("A34").Value = 'HTS_TST.STI.Pos.T'
$ws.Range("A35").Value = 'HTS_RECENT.STI.T'
$ws.Range("A36").Value = 'HTS_TST.VCT.Neg.T'
$ws.Range("A37").Value = 'HTS_TST.VCT.Pos.T'
$ws.Range("A38").Value = 'HTS_RECENT.VCT.T'
$ws.Range("A39").Value = 'HTS_TST.VCT.Neg.T'
$ws.Range("A40").Value = 'HTS_TST.VCT.Pos.T'
$ws.Range("A41").Value = 'HTS_TST.KP.Neg.T'
$ws.Range("A42").Value = 'HTS_TST.KP.Pos.T'
$ws.Range("A43").Value = 'HTS_SELF.KP.T'
$ws.Range("A44").Value = 'HTS_RECENT.KP.T'
$ws.Range("A45").Value = 'KP_MAT.T'
$ws.Range("A46").Value = 'KP_PREV.T'
$ws.Range("A47").Value = 'OVC_HIVSTAT.T'
$ws.Range("A48").Value = 'OVC_SERV.Active.T'
$ws.Range("A49").Value = 'OVC_SERV.Graduated.T'
$ws.Range("A50").Value = 'PMTCT_ART.Already.T'
$ws.Range("A51").Value = 'PMTCT_ART.New.T'
$ws.Range("A52").Value = 'PMTCT_EID.N.2.T'
$ws.Range("A53").Value = 'PMTCT_EID.N.12.T'
$ws.Range("A54").Value = 'PMTCT_STAT.D.T'
$ws.Range("A55").Value = 'PMTCT_STAT.N.KnownPos.T'
$ws.Range("A56").Value = 'PMTCT_STAT.N.New.Neg.T'
$ws.Range("A57").Value = 'PMTCT_STAT.N.New.Pos.T'
$ws.Range("A58").Value = 'HTS_RECENT.PMTCT.T'
$ws.Range("A59").Value = 'PP_PREV.T'
$ws.Range("A60").Value = 'PrEP_CURR.T'
$ws.Range("A61").Value = 'PrEP_CURR.KP.T'
$ws.Range("A62").Value = 'PrEP_NEW.T'
$ws.Range("A63").Value = 'PrEP_NEW.KP.T'
$ws.Range("A64").Value = 'TB_ART.Already.T'
$ws.Range("A65").Value = 'TB_ART.New.T'
$ws.Range("A66").Value = 'TB_PREV.D.Already.T'
$ws.Range("A67").Value = 'TB_PREV.D.New.T'
$ws.Range("A68").Value = 'TB_PREV.N.Already.T'
$ws.Range("A69").Value = 'TB_PREV.N.New.T'
$ws.Range("A70").Value = 'TB_STAT.D.T'
$ws.Range("A71").Value = 'TB_STAT.N.KnownPos.T'
$ws.Range("A72").Value = 'TB_STAT.N.New.Neg.T'
$ws.Range("A73").Value = 'TB_STAT.N.New.Pos.T'
$ws.Range("A74").Value = 'HTS_RECENT.TB.T'
$ws.Range("A75").Value = 'TX_CURR.T'
$ws.Range("A76").Value = 'TX_NEW.T'
$ws.Range("A77").Value = 'TX_NEW.KP.T'
$ws.Range("A78").Value = 'TX_CURR.KP.T'
$ws.Range("A79").Value = 'TX_PVLS.D.KP.T'
$ws.Range("A80").Value = 'TX_PVLS.N.KP.T'
$ws.Range("A81").Value = 'TX_PVLS.D.Routine.T'
$ws.Range("A82").Value = 'TX_PVLS.N.Routine.T'
$ws.Range("A83").Value = 'TX_TB.D.Already.Neg.T'
$ws.Range("A84").Value = 'TX_TB.D.New.Neg.T'
$ws.Range("A85").Value = 'TX_TB.D.Already.Pos.T'
$ws.Range("A86").Value = 'TX_TB.D.New.Pos.T'
$ws.Range("A87").Value = 'VMMC_CIRC.Neg.T'
$ws.Range("A88").Value = 'VMMC_CIRC.Pos.T'
$ws.Range("A89").Value = 'HTS_RECENT.VMMC.T'
$ws.Range("A90").Value = 'VMMC_CIRC.Unk.T'
$ws.Range("A91").Value = 'CXCA_SCRN.T'
